$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update total "VALOR MORA" amount
$ws.Range("E11").Value = 104000

# Update "Cant. Periodos" count
$ws.Range("F13").Value = 2

# The period that used to be "2507" is now "2405", with an updated value
$ws.Range("E16").Value = "2405"
$ws.Range("G16").Value = 1300000

# Update the "2406" row's value as well
$ws.Range("G17").Value = 1300000

# Row 18 carries the "last row" bottom-border formatting; copy that formatting
# onto row 17 before removing row 18, so the new last data row keeps the
# closing border of the table.
$ws.Range("B18:J18").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Remove the old "2405" row (row 18), which is now redundant after consolidating
# into row 16; this shifts the remaining rows (legal representative block) up by one
$ws.Rows("18:18").Delete()
